$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor data corrections
$ws.Range("E5").Value = 6372
$ws.Range("E6").Value = 29968
$ws.Range("E20").Value = 5630835

# Update selection to E7
$ws.Range("E7").Select()
